$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 50, shifting existing rows (50-91) down to (51-92)
$ws.Rows.Item(50).Insert()

# Populate the new row 50 with its values.
# Columns A,B,C,E,F,G,H,I,J,K,L,Q,T mirror the row that used to be at 50 (now 51)
$ws.Range("A50").Value = 10
$ws.Range("B50").Value = "Vega Modelo de Temuco"
$ws.Range("C50").Value = "La Araucanía"
$ws.Range("D50").Value = 44669
$ws.Range("E50").Value = 9
$ws.Range("F50").Value = "Fruta"
$ws.Range("G50").Value = 100101
$ws.Range("H50").Value = "Berries"
$ws.Range("I50").Value = 100101001
$ws.Range("J50").Value = "Arándano (blue)"
$ws.Range("K50").Value = "Sin especificar"
$ws.Range("L50").Value = "Primera"
$ws.Range("M50").Value = 100
$ws.Range("N50").Value = 3000
$ws.Range("O50").Value = 3000
$ws.Range("P50").Value = 3000
$ws.Range("Q50").Value = "$/kilo"
$ws.Range("R50").Value = "Región de La Araucanía"
$ws.Range("S50").Value = 3000
$ws.Range("T50").Value = 1
